$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D (Price) and E (Volume) columns retain their original text (string) storage,
# matching how the source data already stores these as inline strings rather than numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.167.43'
$ws.Range('E2').Value = '  -0.70%  '
$ws.Range('D3').Value = '1.846.39'
$ws.Range('E3').Value = '  -0.80%  '
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '313.65'
$ws.Range('E5').Value = '  -0.46%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D7').Value = '0.4634'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.3694'
$ws.Range('E8').Value = '  -0.89%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').Value = '0.8856'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('D11').Value = '19.89'
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').Value = '0.07819'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').Value = '1.860.98'
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').Value = '5.395'
$ws.Range('E14').Value = '  -0.50%  '
$ws.Range('D15').Value = '6.503'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').Value = '91.41'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '1.004'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '0.000008839'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = '1.003'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '27.200.80'
$ws.Range('E20').Value = '  -0.53%  '
$ws.Range('D21').Value = '14.62'
$ws.Range('E21').Value = '  -2.04%  '
$ws.Range('D22').Value = '5.040'
$ws.Range('E22').Value = '  -2.23%  '
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').Value = '2.030'
$ws.Range('E24').Value = '  +8.75%  '
$ws.Range('D25').Value = '150.81'
$ws.Range('E25').Value = '  -1.53%  '
$ws.Range('D26').Value = '18.36'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('D27').Value = '2.024'
$ws.Range('E27').Value = '  -2.96%  '
$ws.Range('D28').Value = '115.66'
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('D29').Value = '5.008'
$ws.Range('E29').Value = '  -2.95%  '
$ws.Range('D30').Value = '0.08839'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('B31').Value = 'HuobiToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D31').Value = '3.157'
$ws.Range('E31').Value = '  +6.00%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '0.7814'
$ws.Range('E32').Value = '  +3.72%  '
$ws.Range('D33').Value = '4.502'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').Value = '1.149'
$ws.Range('E34').Value = '  -0.36%  '
$ws.Range('D35').Value = '2.690'
$ws.Range('E35').Value = '  +3.33%  '
$ws.Range('D36').Value = '1.102'
$ws.Range('E36').Value = '  +1.82%  '
$ws.Range('D37').Value = '0.01943'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').Value = '0.05210'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').Value = '2.950'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').Value = '7.006'
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('D41').Value = '0.5038'
$ws.Range('E41').Value = '  -3.25%  '
$ws.Range('D42').Value = '0.1610'
$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('D43').Value = '8.486'
$ws.Range('E43').Value = '  +1.59%  '
$ws.Range('D44').Value = '0.4744'
$ws.Range('E44').Value = '  -3.21%  '
$ws.Range('D45').Value = '10.31'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('D46').Value = '1.002'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').Value = '102.98'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('D48').Value = '1.630'
$ws.Range('E48').Value = '  -0.66%  '
$ws.Range('D49').Value = '0.06189'
$ws.Range('E49').Value = '  -1.39%  '
$ws.Range('D50').Value = '65.38'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').Value = '36.44'
$ws.Range('E51').Value = '  -2.00%  '
